$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("M2").Value = "[49.92080998514051, 50.15285148640196]"
$ws.Range("U2").Value = "[49.933974043296026, 50.09632974622439]"
$ws.Range("M3").Value = "[49.87814390862527, 50.1272732849089]"
$ws.Range("U3").Value = "[49.997376600394155, 50.153590264010354]"
$ws.Range("M4").Value = "[49.827221175097996, 50.1618703901272]"
$ws.Range("U4").Value = "[49.93477441921237, 50.1241546897752]"
$ws.Range("M5").Value = "[49.94785563942253, 50.25116474286857]"
$ws.Range("U5").Value = "[49.95245053225074, 50.11474534756101]"
$ws.Range("M6").Value = "[49.82482253319466, 50.10540886136671]"
$ws.Range("U6").Value = "[49.99531956312543, 50.13982951753683]"
$ws.Range("M7").Value = "[49.90930403224561, 50.194528316759495]"
$ws.Range("U7").Value = "[50.011773336001205, 50.162337531628815]"
$ws.Range("M8").Value = "[49.849953881211285, 50.187755720028996]"
$ws.Range("U8").Value = "[49.93129528888244, 50.11242080809639]"
$ws.Range("M9").Value = "[49.856932177566, 50.17319537740726]"
$ws.Range("U9").Value = "[49.95041495549029, 50.12048485359087]"
$ws.Range("M10").Value = "[49.79315278222805, 50.084721897802254]"
$ws.Range("U10").Value = "[49.91676981906598, 50.08480500035447]"
$ws.Range("M11").Value = "[49.898985856636656, 50.219723676526186]"
$ws.Range("U11").Value = "[50.01243949001261, 50.20061956794741]"
$ws.Range("M12").Value = "[49.94708497079242, 50.141262465083116]"
$ws.Range("U12").Value = "[49.965739263474546, 50.10463138031213]"
$ws.Range("M13").Value = "[49.88521646297608, 50.18639019483897]"
$ws.Range("U13").Value = "[49.87308503388744, 50.04518015581491]"
$ws.Range("M14").Value = "[49.858860136961745, 50.17573640966064]"
$ws.Range("U14").Value = "[49.987372999723846, 50.154909687745956]"
$ws.Range("M15").Value = "[49.769694449018004, 50.09239488835709]"
$ws.Range("U15").Value = "[49.94871041264968, 50.117452415862324]"
$ws.Range("M16").Value = "[49.81847867059839, 50.109656032066205]"
$ws.Range("U16").Value = "[49.9173694429299, 50.09126184942573]"
$ws.Range("M17").Value = "[49.86831495434458, 50.13675900270504]"
$ws.Range("U17").Value = "[49.87802088701233, 50.04832398023685]"
$ws.Range("M18").Value = "[49.82163330362212, 50.13754040334172]"
$ws.Range("U18").Value = "[49.92319987416617, 50.09406621294279]"
